# Applies the commit:
#   fix: ensure elective courses are scheduled in same time slots for both
#   sections A and B
#
# 1) Renumber/reshuffle courses in the Section_A and Section_B timetables
#    (course codes CS312/CS307/CS308/CS465 become CS309/CS303/CS304/CS461,
#    plus the elective CS461 is realigned to the same Thu 15:30-17:00 slot
#    in both sections).
# 2) Update the Course_Summary sheet with the renumbered course codes, new
#    names and new instructors.
# 3) Add a new Elective_Coordination sheet documenting the shared elective
#    slot for sections A & B.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Section_A timetable updates
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")
$wsA.Range("C2").Value = "CS304"
$wsA.Range("D2").Value = "CS304"
$wsA.Range("F2").Value = "CS304"
$wsA.Range("B3").Value = "CS309 (Tutorial)"
$wsA.Range("E3").Value = "Free"
$wsA.Range("F3").Value = "Free"
$wsA.Range("C5").Value = "CS303"
$wsA.Range("D5").Value = "Free"
$wsA.Range("E5").Value = "CS461 (Elective)"
$wsA.Range("F5").Value = "Free"
$wsA.Range("B6").Value = "CS309"
$wsA.Range("C6").Value = "CS309"
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "CS309"
$wsA.Range("F6").Value = "CS303 (Tutorial)"
$wsA.Range("C7").Value = "CS304 (Tutorial)"
$wsA.Range("D7").Value = "CS303"
$wsA.Range("F7").Value = "CS303"

# ---------------------------------------------------------------------
# Section_B timetable updates
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")
$wsB.Range("B2").Value = "CS303"
$wsB.Range("D2").Value = "CS304 (Tutorial)"
$wsB.Range("E2").Value = "Free"
$wsB.Range("F2").Value = "Free"
$wsB.Range("B3").Value = "Free"
$wsB.Range("D3").Value = "CS303 (Tutorial)"
$wsB.Range("E3").Value = "CS304"
$wsB.Range("B5").Value = "CS309 (Tutorial)"
$wsB.Range("C5").Value = "CS304"
$wsB.Range("E5").Value = "CS461 (Elective)"
$wsB.Range("C6").Value = "CS303"
$wsB.Range("D6").Value = "CS304"
$wsB.Range("E6").Value = "CS303"
$wsB.Range("B7").Value = "Free"
$wsB.Range("C7").Value = "CS309"
$wsB.Range("D7").Value = "CS309"
$wsB.Range("E7").Value = "CS309"

# ---------------------------------------------------------------------
# Course_Summary updates
# ---------------------------------------------------------------------
$wsC = $wb.Worksheets.Item("Course_Summary")
$wsC.Range("A2").Value = "CS309"
$wsC.Range("B2").Value = "Statistics for Computer Science"
$wsC.Range("F2").Value = "Dr. Sunil C K"
$wsC.Range("A3").Value = "CS303"
$wsC.Range("B3").Value = "Computer Networks"
$wsC.Range("F3").Value = "Dr. Animesh Roy"
$wsC.Range("A4").Value = "CS304"
$wsC.Range("B4").Value = "Artificial Intelligence"
$wsC.Range("F4").Value = "Dr. Krishnendu"
$wsC.Range("A5").Value = "CS461"
$wsC.Range("B5").Value = "Parallel computing"
$wsC.Range("F5").Value = "Dr. Pramod"

# ---------------------------------------------------------------------
# New Elective_Coordination sheet (added after Course_Summary, i.e. last)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsE = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsE.Name = "Elective_Coordination"

# Reuse the same header formatting (bold, bordered, centered) already used
# by the header rows on the other sheets.
$wsC.Range("A1:E1").Copy()
$wsE.Range("A1:E1").PasteSpecial(-4122)

$wsE.Cells.Item(1, 1).Value = "Elective Course"
$wsE.Cells.Item(1, 2).Value = "Day"
$wsE.Cells.Item(1, 3).Value = "Time Slot"
$wsE.Cells.Item(1, 4).Value = "Slot Name"
$wsE.Cells.Item(1, 5).Value = "Sections"

$wsE.Cells.Item(2, 1).Value = "CS461"
$wsE.Cells.Item(2, 2).Value = "Thu"
$wsE.Cells.Item(2, 3).Value = "15:30-17:00"
$wsE.Cells.Item(2, 4).Value = "Elective_Slot_1"
$wsE.Cells.Item(2, 5).Value = "A & B (Common Slot)"
